$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for TAS (Tashkent, Uzbekistan) at row 235.
# This shifts all subsequent rows up by one and shrinks the used range.
$ws.Rows.Item(235).Delete()
